$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Fecha/Volumen/Precios values between row 2 <-> row 4, and row 3 <-> row 7.
# Columns involved: D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), P (Precio $/Kg)
$cols = @("D", "J", "K", "L", "M", "P")

function Swap-RowValues($ws, $row1, $row2, $cols) {
    foreach ($col in $cols) {
        $addr1 = "$col$row1"
        $addr2 = "$col$row2"
        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2
        $ws.Range($addr1).Value = $v2
        $ws.Range($addr2).Value = $v1
    }
}

Swap-RowValues $ws 2 4 $cols
Swap-RowValues $ws 3 7 $cols
